$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5196889042854309
$ws.Range("B1").Value = 2.399427652359009
$ws.Range("C1").Value = 6.187565326690674
$ws.Range("D1").Value = 1.565920829772949
$ws.Range("E1").Value = 0.900866687297821
